$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header key row (row 1) with the new field keys
$ws.Range("A1").Value = "peminatancode"
$ws.Range("B1").Value = "peminatanname"

# Change the active selection to C1 (as reflected in the saved file)
$ws.Range("C1").Select()

$wb.Save()
